$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 303.69232
$ws.Range("I53").Value = 217.41176
$ws.Range("J53").Value = 466.66666
$ws.Range("K53").Value = 217.41176
$ws.Range("L53").Value = 466.66666
$ws.Range("M53").Value = 419.58824
$ws.Range("N53").Value = -1740.66666
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 20000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -21248
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 60000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -66240
$ws.Range("H98").Value = 623064.9
$ws.Range("I98").Value = 1017531.44
$ws.Range("J98").Value = 3188.8572
$ws.Range("K98").Value = 1017531.44
$ws.Range("L98").Value = 3188.8572
$ws.Range("M98").Value = -1016033.44
$ws.Range("N98").Value = -6184.8572
$ws.Range("H100").Value = 9261194
$ws.Range("I100").Value = 12822568
$ws.Range("J100").Value = 1622.4
$ws.Range("K100").Value = 12822568
$ws.Range("L100").Value = 1622.4
$ws.Range("M100").Value = -12822027
$ws.Range("N100").Value = -2704.4
$ws.Range("H107").Value = 483488.44
$ws.Range("I107").Value = 505419.72
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 505419.72
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -503499.72
$ws.Range("N107").Value = -4840
$ws.Range("H122").Value = 623064.9
$ws.Range("I122").Value = 1017531.44
$ws.Range("J122").Value = 3188.8572
$ws.Range("K122").Value = 3052594.32
$ws.Range("L122").Value = 9566.571599999999
$ws.Range("M122").Value = -3050144.32
$ws.Range("N122").Value = -14466.5716
$ws.Range("H133").Value = 37666.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 37666.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 37666.5
$ws.Range("N133").Value = -47786.5
$ws.Range("H137").Value = 16129955
$ws.Range("I137").Value = 19231450
$ws.Range("J137").Value = 2177.8
$ws.Range("K137").Value = 57694350
$ws.Range("L137").Value = 6533.400000000001
$ws.Range("M137").Value = -57691800
$ws.Range("N137").Value = -11633.4

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22413.51
$ws.Range("I32").Value = 3361.4468
$ws.Range("J32").Value = 171654.67
$ws.Range("K32").Value = 3361.4468
$ws.Range("L32").Value = 171654.67
$ws.Range("M32").Value = -3074.4468
$ws.Range("N32").Value = -172228.67
$ws.Range("H74").Value = 2905.3428
$ws.Range("I74").Value = 884.04083
$ws.Range("J74").Value = 7621.7144
$ws.Range("K74").Value = 884.04083
$ws.Range("L74").Value = 7621.7144
$ws.Range("M74").Value = -10.04083000000003
$ws.Range("N74").Value = -9369.714400000001
$ws.Range("H77").Value = 2905.3428
$ws.Range("I77").Value = 884.04083
$ws.Range("J77").Value = 7621.7144
$ws.Range("K77").Value = 4420.20415
$ws.Range("L77").Value = 38108.572
$ws.Range("M77").Value = -52.20415000000048
$ws.Range("N77").Value = -46844.572
$ws.Range("H102").Value = 2513.3333
$ws.Range("I102").Value = 2577.5
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2577.5
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -955.5
$ws.Range("N102").Value = -5244
$ws.Range("H133").Value = 50966.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50966.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50966.332
$ws.Range("N133").Value = -56026.332
$ws.Range("H139").Value = 40261.75
$ws.Range("I139").Value = 30470
$ws.Range("J139").Value = 43525.668
$ws.Range("K139").Value = 30470
$ws.Range("L139").Value = 43525.668
$ws.Range("M139").Value = -25330
$ws.Range("N139").Value = -53805.668
$ws.Range("H141").Value = 44999.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 44999.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 44999.5
$ws.Range("N141").Value = -55359.5

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1278.4706
$ws.Range("I94").Value = 1103.8572
$ws.Range("J94").Value = 2093.3333
$ws.Range("K94").Value = 1103.8572
$ws.Range("L94").Value = 2093.3333
$ws.Range("M94").Value = -652.8571999999999
$ws.Range("N94").Value = -2995.3333
$ws.Range("H105").Value = 3427.4814
$ws.Range("I105").Value = 3344.3125
$ws.Range("J105").Value = 3548.4546
$ws.Range("K105").Value = 3344.3125
$ws.Range("L105").Value = 3548.4546
$ws.Range("M105").Value = -1597.3125
$ws.Range("N105").Value = -7042.4546
$ws.Range("H120").Value = 35000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 35000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676
$ws.Range("H133").Value = 50000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120
$ws.Range("H134").Value = 15626766
$ws.Range("I134").Value = 19232138
$ws.Range("J134").Value = 3493.4167
$ws.Range("K134").Value = 57696414
$ws.Range("L134").Value = 10480.2501
$ws.Range("M134").Value = -57693879
$ws.Range("N134").Value = -15550.2501

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 540.0476
$ws.Range("I22").Value = 489.46155
$ws.Range("J22").Value = 622.25
$ws.Range("K22").Value = 489.46155
$ws.Range("L22").Value = 622.25
$ws.Range("M22").Value = -139.46155
$ws.Range("N22").Value = -1322.25
$ws.Range("H31").Value = 1239.7937
$ws.Range("I31").Value = 802.50946
$ws.Range("J31").Value = 3557.4
$ws.Range("K31").Value = 802.50946
$ws.Range("L31").Value = 3557.4
$ws.Range("M31").Value = -507.50946
$ws.Range("N31").Value = -4147.4
$ws.Range("H34").Value = 1239.7937
$ws.Range("I34").Value = 802.50946
$ws.Range("J34").Value = 3557.4
$ws.Range("K34").Value = 802.50946
$ws.Range("L34").Value = 3557.4
$ws.Range("M34").Value = -600.50946
$ws.Range("N34").Value = -3961.4
$ws.Range("H64").Value = 35000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 35000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
$ws.Range("H67").Value = 35000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 35000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716
$ws.Range("H99").Value = 15625875
$ws.Range("I99").Value = 15625875
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 15625875
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -15624377
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 861.8570999999999
$ws.Range("I105").Value = 867.46155
$ws.Range("J105").Value = 789
$ws.Range("K105").Value = 867.46155
$ws.Range("L105").Value = 789
$ws.Range("M105").Value = 879.53845
$ws.Range("N105").Value = -4283
$ws.Range("H126").Value = 15625875
$ws.Range("I126").Value = 15625875
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 46877625
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -46875155
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2193.0815
$ws.Range("I132").Value = 1612.925
$ws.Range("J132").Value = 4771.5557
$ws.Range("K132").Value = 4838.775
$ws.Range("L132").Value = 14314.6671
$ws.Range("M132").Value = -2308.775
$ws.Range("N132").Value = -19374.6671

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8986.441000000001
$ws.Range("I4").Value = 111.117645
$ws.Range("J4").Value = 17861.766
$ws.Range("K4").Value = 333.352935
$ws.Range("L4").Value = 53585.298
$ws.Range("M4").Value = -221.352935
$ws.Range("N4").Value = -53809.298
$ws.Range("H113").Value = 16667635
$ws.Range("I113").Value = 644.7143
$ws.Range("J113").Value = 31251252
$ws.Range("K113").Value = 1934.1429
$ws.Range("L113").Value = 93753756
$ws.Range("M113").Value = 235.8571000000002
$ws.Range("N113").Value = -93758096
$ws.Range("H131").Value = 1658.9656
$ws.Range("I131").Value = 320
$ws.Range("J131").Value = 2169.0476
$ws.Range("K131").Value = 960
$ws.Range("L131").Value = 6507.1428
$ws.Range("M131").Value = 4080
$ws.Range("N131").Value = -16587.1428

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 20000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 20000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496
$ws.Range("H67").Value = 20000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 20000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716
$ws.Range("H122").Value = 1112349.8
$ws.Range("I122").Value = 1390062.1
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4170186.3
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -4167736.3
$ws.Range("N122").Value = -9400
$ws.Range("H138").Value = 64100
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 64100
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 64100
$ws.Range("N138").Value = -74380
$ws.Range("H139").Value = 35000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 35000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2987.0476
$ws.Range("I7").Value = 1954
$ws.Range("J7").Value = 3503.5715
$ws.Range("K7").Value = 1954
$ws.Range("L7").Value = 3503.5715
$ws.Range("M7").Value = -1842
$ws.Range("N7").Value = -3727.5715
$ws.Range("H40").Value = 3159.5925
$ws.Range("I40").Value = 1225.4
$ws.Range("J40").Value = 4297.353
$ws.Range("K40").Value = 1225.4
$ws.Range("L40").Value = 4297.353
$ws.Range("M40").Value = -1089.4
$ws.Range("N40").Value = -4569.353
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H126").Value = 2987.0476
$ws.Range("I126").Value = 1954
$ws.Range("J126").Value = 3503.5715
$ws.Range("K126").Value = 5862
$ws.Range("L126").Value = 10510.7145
$ws.Range("M126").Value = -3392
$ws.Range("N126").Value = -15450.7145

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 21000.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 21000.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 21000.5
$ws.Range("N80").Value = -22996.5
$ws.Range("H83").Value = 21000.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 21000.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 63001.5
$ws.Range("N83").Value = -72985.5
$ws.Range("H107").Value = 5556806.5
$ws.Range("I107").Value = 9260181
$ws.Range("J107").Value = 1745
$ws.Range("K107").Value = 27780543
$ws.Range("L107").Value = 5235
$ws.Range("M107").Value = -27778623
$ws.Range("N107").Value = -9075
$ws.Range("H132").Value = 11113465
$ws.Range("I132").Value = 16669035
$ws.Range("J132").Value = 2323.6
$ws.Range("K132").Value = 50007105
$ws.Range("L132").Value = 6970.799999999999
$ws.Range("M132").Value = -50004575
$ws.Range("N132").Value = -12030.8
